# Generate Report for Handoff
# A new handoff round completed for the "431899e4-..." file: refresh the
# "Latest Handoff Datetime" (column D) on both locale report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-11 04:58:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-11 04:58:21"
